$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("A50").Value = 'Start dag'
$ws.Range("B50").Value = 0
$ws.Range("C49").Copy($ws.Range("C50"))
$ws.Range("C50").Value = 44692
$ws.Range("K50").Formula = '=IF(C50<>0,DAYS360($F$1,C50),0)'
$ws.Range("L50").Formula = '=L49-B50'

$ws.Range("A51").Value = 'Indkøbskurv'
$ws.Range("B51").Value = 0.25
$ws.Range("C49").Copy($ws.Range("C51"))
$ws.Range("C51").Value = 44692
$ws.Range("K51").Formula = '=IF(C51<>0,DAYS360($F$1,C51),0)'
$ws.Range("L51").Formula = '=L50-B51'

$ws.Range("A52").Value = 'UI Check ud '
$ws.Range("B52").Value = 6
$ws.Range("C49").Copy($ws.Range("C52"))
$ws.Range("C52").Value = 44692
$ws.Range("K52").Formula = '=IF(C52<>0,DAYS360($F$1,C52),0)'
$ws.Range("L52").Formula = '=L51-B52'

$ws.Range("A53").Value = 'Settings'
$ws.Range("B53").Value = 0.55
$ws.Range("C49").Copy($ws.Range("C53"))
$ws.Range("C53").Value = 44692
$ws.Range("K53").Formula = '=IF(C53<>0,DAYS360($F$1,C53),0)'
$ws.Range("L53").Formula = '=L52-B53'

$ws.Range("A54").Value = 'Vare information'
$ws.Range("B54").Value = 2.05
$ws.Range("C49").Copy($ws.Range("C54"))
$ws.Range("C54").Value = 44692
$ws.Range("K54").Formula = '=IF(C54<>0,DAYS360($F$1,C54),0)'
$ws.Range("L54").Formula = '=L53-B54'

$ws.Range("A55").Value = 'Start dag'
$ws.Range("B55").Value = 0
$ws.Range("C49").Copy($ws.Range("C55"))
$ws.Range("C55").Value = 44693
$ws.Range("K55").Formula = '=IF(C55<>0,DAYS360($F$1,C55),0)'
$ws.Range("L55").Formula = '=L54-B55'

$ws.Range("A56").Value = 'Database wallet'
$ws.Range("B56").Value = 0.26
$ws.Range("C49").Copy($ws.Range("C56"))
$ws.Range("C56").Value = 44693
$ws.Range("K56").Formula = '=IF(C56<>0,DAYS360($F$1,C56),0)'
$ws.Range("L56").Formula = '=L55-B56'

$ws.Range("A57").Value = 'Database check ud'
$ws.Range("B57").Value = 0.6
$ws.Range("C49").Copy($ws.Range("C57"))
$ws.Range("C57").Value = 44693
$ws.Range("K57").Formula = '=IF(C57<>0,DAYS360($F$1,C57),0)'
$ws.Range("L57").Formula = '=L56-B57'

$ws.Range("A58").Value = 'Opret kort'
$ws.Range("B58").Value = 0.16
$ws.Range("C49").Copy($ws.Range("C58"))
$ws.Range("C58").Value = 44693
$ws.Range("K58").Formula = '=IF(C58<>0,DAYS360($F$1,C58),0)'
$ws.Range("L58").Formula = '=L57-B58'

$ws.Range("A59").Value = 'Start dag'
$ws.Range("B59").Value = 0
$ws.Range("C49").Copy($ws.Range("C59"))
$ws.Range("C59").Value = 44694
$ws.Range("K59").Formula = '=IF(C59<>0,DAYS360($F$1,C59),0)'
$ws.Range("L59").Formula = '=L58-B59'

$ws.Range("A60").Value = 'Gemte kurv'
$ws.Range("B60").Value = 0.26
$ws.Range("C49").Copy($ws.Range("C60"))
$ws.Range("C60").Value = 44694
$ws.Range("K60").Formula = '=IF(C60<>0,DAYS360($F$1,C60),0)'
$ws.Range("L60").Formula = '=L59-B60'

$ws.Range("A61").Value = 'Database gemte kurv'
$ws.Range("B61").Value = 0.26
$ws.Range("C49").Copy($ws.Range("C61"))
$ws.Range("C61").Value = 44694
$ws.Range("K61").Formula = '=IF(C61<>0,DAYS360($F$1,C61),0)'
$ws.Range("L61").Formula = '=L60-B61'

$ws.Range("A62").Value = 'Check ud'
$ws.Range("B62").Value = 0.6
$ws.Range("C49").Copy($ws.Range("C62"))
$ws.Range("C62").Value = 44694
$ws.Range("K62").Formula = '=IF(C62<>0,DAYS360($F$1,C62),0)'
$ws.Range("L62").Formula = '=L61-B62'

$ws.Range("A63").Value = 'Wallet'
$ws.Range("B63").Value = 0.6
$ws.Range("C49").Copy($ws.Range("C63"))
$ws.Range("C63").Value = 44694
$ws.Range("K63").Formula = '=IF(C63<>0,DAYS360($F$1,C63),0)'
$ws.Range("L63").Formula = '=L62-B63'

$ws.Range("A64").Value = 'Start dag'
$ws.Range("B64").Value = 0
$ws.Range("C49").Copy($ws.Range("C64"))
$ws.Range("C64").Value = 44695
$ws.Range("K64").Formula = '=IF(C64<>0,DAYS360($F$1,C64),0)'
$ws.Range("L64").Formula = '=L63-B64'

$ws.Range("A65").Value = 'Dataindsættelse indkøbskurv'
$ws.Range("B65").Value = 20.6
$ws.Range("C49").Copy($ws.Range("C65"))
$ws.Range("C65").Value = 44695
$ws.Range("K65").Formula = '=IF(C65<>0,DAYS360($F$1,C65),0)'
$ws.Range("L65").Formula = '=L64-B65'

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("M59").Select()

Write-Host "edit complete"
